$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.344089388847351
$ws.Range("B1").Value = 2.514552116394043
$ws.Range("C1").Value = 5.771359443664551
$ws.Range("D1").Value = 1.917012333869934
$ws.Range("E1").Value = 1.255420088768005
